$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 51 also has Coin name and Link changes
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"

$ws.Range("D2").Value = "57.961.37"
$ws.Range("E2").Value = "  -0.75%  "

$ws.Range("D3").Value = "2.574.92"
$ws.Range("E3").Value = "  -2.74%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "517.88"
$ws.Range("E5").Value = "  -0.66%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.97"
$ws.Range("E6").Value = "  -1.26%  "

$ws.Range("E7").Value = "  -0.21%  "

$ws.Range("E8").Value = "  -1.26%  "

$ws.Range("D9").Value = "2.590.14"
$ws.Range("E9").Value = "  -2.34%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.58"
$ws.Range("E10").Value = "  -2.05%  "

$ws.Range("E11").Value = "  -1.66%  "

$ws.Range("E12").Value = "  -4.67%  "

$ws.Range("E13").Value = "  -1.00%  "

$ws.Range("D14").Value = "3.028.71"
$ws.Range("E14").Value = "  -2.83%  "

$ws.Range("D15").Value = "57.936.89"
$ws.Range("E15").Value = "  -0.82%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.28"
$ws.Range("E16").Value = "  -2.92%  "

$ws.Range("D18").Value = "2.560.58"
$ws.Range("E18").Value = "  -3.62%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "340.00"
$ws.Range("E19").Value = "  +0.37%  "

$ws.Range("E20").Value = "  -2.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.20"
$ws.Range("E21").Value = "  -2.15%  "

$ws.Range("E22").Value = "  -0.09%  "

$ws.Range("E23").Value = "  -0.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.41"
$ws.Range("E24").Value = "  +1.56%  "

$ws.Range("E25").Value = "  -0.42%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.401"
$ws.Range("E26").Value = "  -5.77%  "

$ws.Range("E27").Value = "  -0.32%  "

$ws.Range("D28").Value = "2.686.65"
$ws.Range("E28").Value = "  -3.01%  "

$ws.Range("E29").Value = "  -2.41%  "

$ws.Range("D30").Value = "0.0₃0747"
$ws.Range("E30").Value = "  -6.32%  "

$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("E32").Value = "  -5.43%  "

$ws.Range("E33").Value = "  -0.77%  "

$ws.Range("E34").Value = "  -1.32%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "149.58"
$ws.Range("E35").Value = "  -1.69%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.00"
$ws.Range("E36").Value = "  -3.87%  "

$ws.Range("E37").Value = "  -3.42%  "

$ws.Range("E38").Value = "  -4.51%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.13"
$ws.Range("E39").Value = "  -1.89%  "

$ws.Range("E40").Value = "  +1.30%  "

$ws.Range("E41").Value = "  -4.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.52"
$ws.Range("E42").Value = "  -2.65%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.997"
$ws.Range("E43").Value = "  -0.32%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "270.55"
$ws.Range("E44").Value = "  -1.30%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.65"
$ws.Range("E45").Value = "  +0.27%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0953"
$ws.Range("E46").Value = "  -2.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.587"
$ws.Range("E47").Value = "  -3.60%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.80"
$ws.Range("E48").Value = "  -3.36%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0521"
$ws.Range("E49").Value = "  -2.44%  "

$ws.Range("D50").Value = "1.976.93"
$ws.Range("E50").Value = "  -3.22%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.60"
$ws.Range("E51").Value = "  -1.93%  "
